$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A7").Value = 43985
$ws.Range("B7").Value = "6 hours 24 minutes"
$ws.Range("A13").Value = "TOTAL:"

$ws.Range("B13").Select()
